$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("values")

# Update PANORAMA_IP value
$ws.Range("B4").Value = "192.168.55.8"

# Remove the INTERNET_ZONE row (row 24) entirely, shifting rows below up
$ws.Rows.Item(24).Delete()
